# File & Task INDEX workbook update:
#  - add a "venue checklist" document row to General_Information section
#  - add a "venue correspondence email" row to Email_Templates section
#  - rename "Poster judging preliminary results" -> "Poster judging template"
#  - rename "All submitted talk abstracts and voting" -> "All submitted talk abstracts"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-FilePathCell {
    param(
        $cell,
        [string]$folder,
        [string]$file
    )
    $full = $folder + $file
    $cell.Value = $full
    $folderLen = $folder.Length
    $fileLen = $file.Length

    $folderChars = $cell.Characters(1, $folderLen)
    $folderChars.Font.Name = "Calibri (Body)"
    $folderChars.Font.Size = 12
    $folderChars.Font.Color = 10921638

    $fileChars = $cell.Characters($folderLen + 1, $fileLen)
    $fileChars.Font.Name = "Calibri"
    $fileChars.Font.Size = 12
    $fileChars.Font.Color = 0
}

# 1) Insert new row in the General_Information block (between
#    "Final breakdown of budget" and "Giant check presented...")
$ws.Rows("101").Insert()
$ws.Range("B101").Value = "Checklist of items for venue and food"
Set-FilePathCell -cell $ws.Range("C101") -folder "2022_Symposium/General_Information/" -file "Generic_Venue_Checklist_060722.docx"

# 2) Insert new row in the Email_Templates block (between
#    "Deadline reminder to register for the symposium" and "Keynote invitation to speak")
$ws.Rows("92").Insert()
$ws.Range("B92").Value = "Email to venue to initate planning"
Set-FilePathCell -cell $ws.Range("C92") -folder "2022_Symposium/Email_Templates/" -file "Generic_Initial_Venue_Correspondence_060722.docx"

# 3) Rename existing entries (rows have shifted down by 2 due to the inserts above)
$ws.Range("B134").Value = "Poster judging template"
$ws.Range("B140").Value = "All submitted talk abstracts"

# 4) Update view/selection to match the author's final cursor position
$ws.Activate()
$ws.Range("C99").Select()
try {
    $excel.ActiveWindow.ScrollRow = 81
    $excel.ActiveWindow.ScrollColumn = 2
} catch {
}
